$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Examples part 1")

# OpenTbs 1.8.1 beta: rename the "xlsx*" ope keywords to the new common
# "tbs:*" keywords shared between ODS and XLSX templates.

# "Merging data with cell" example
$ws1.Range("C26").Value = "[cell2.score;block=tbs:cell;ope=tbs:num]"

# "Change the type data in a cell" example table
$ws1.Range("C34").Value = "tbs:num"
$ws1.Range("C35").Value = "tbs:bool"
$ws1.Range("C36").Value = "tbs:date"
$ws1.Range("D36").Value = "[onshow.x_dt;ope=tbs:date]"
$ws1.Range("D35").Value = "[onshow.x_bt;ope=tbs:bool]"
$ws1.Range("D34").Value = "[onshow.x_num;ope=tbs:num]"

# "Merging data with rows" example (score column)
$ws1.Range("E20").Value = "[a.score;ope=tbs:num]"
$ws1.Range("F20").Value = "[a.score;ope=tbs:num]"

# New named cell on the "Delete me" sheet, and a workbook-level defined
# name pointing to it.
$ws4 = $wb.Worksheets.Item("Delete me")
$ws4.Range("B6").Value = "And this named cell too."

# New header label for the extra column demonstrating the named cell.
$ws1.Range("F19").Value = "Score again"

$wb.Names.Add("the_named_cell", "='Delete me'!`$B`$6")
